# Advent of Code 2015 - Day 10 ("Elves Look, Elves Say")
# Wire up the Workings sheet: puzzle input in B2, and an array formula in
# B3 that splits the input into its individual digits (B3:K3), plus fix
# the broken "input" defined name so it points at the input cell again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workings")

# The "input" named range had gone stale (#REF!) - point it back at B2.
$wb.Names.Item("input").RefersTo = "=Workings!`$B`$2"

# Puzzle input.
$ws.Range("B2").Value = 1113222113

# Split the input into one digit per cell, spilling from B3 across to K3.
$ws.Range("B3:K3").FormulaArray = "=MID(B2,SEQUENCE(,LEN(B2)),1)"

# Cosmetic touch-ups matching the saved view state.
$ws.Columns.Item(1).ColumnWidth = 11.1640625

[void]$ws.Select()
$excel.ActiveWindow.Zoom = 217
[void]$ws.Range("C4").Select()
